$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New asset rows being documented:
#   Row 6 (existing ConsolidatedFile) now points at the .xlsm (was .xlsx)
#   Row 7 (new) FaultyRowsFile
#   Row 8 (new) Macros
# ---------------------------------------------------------------------------

$processedTarget   = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Processed'
$loadedTarget       = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Loaded'
$templateTarget     = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Template\02%20February%202019-DiCarlo%20Distributors%20Template.xlsx'
$errorTarget        = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Errors'
$consolidatedTarget = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\PriceAudit_Consolidated.xlsm'
$faultyRowsTarget   = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Errors\FaultyRows.xlsm'
$macrosTarget       = 'file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\Macros'

$consolidatedDisplay = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\PriceAudit_Consolidated.xlsm"
$faultyRowsDisplay   = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Errors\FaultyRows.xlsm"
$macrosDisplay       = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\Macros"

# New row content (Name / Description / Asset), written in the same order
# the workbook author entered it so shared-string indices line up.
$ws.Range("A7").Value() = "FaultyRowsFile"
$ws.Range("C7").Value() = "Faulty rows from files with partially correct information"
$ws.Range("B7").Value() = $faultyRowsDisplay
$ws.Range("B6").Value() = $consolidatedDisplay
$ws.Range("A8").Value() = "Macros"
$ws.Range("B8").Value() = $macrosDisplay

# The host only supports clearing hyperlinks for an entire sheet at once, so
# rebuild the full Hyperlinks collection: keep B2:B5 pointed at their
# original targets, then (re)attach B6:B8 with the new targets.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), $processedTarget) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $loadedTarget) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), $templateTarget) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), $errorTarget) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), $consolidatedTarget) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), $faultyRowsTarget) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), $macrosTarget) | Out-Null

# Re-apply the shared "Hyperlink" cell style to every asset link cell.
$ws.Range("B2:B8").Style = "Hyperlink"

# Match the author's final selection / scroll state.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("B8").Select() | Out-Null
